$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Value = "'" + '59.218.46'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(2, 5)
$cell.Value = "'" + '  +1.68%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 4)
$cell.Value = "'" + '2.587.92'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 5)
$cell.Value = "'" + '  -0.34%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(4, 5)
$cell.Value = "'" + '  -0.04%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'" + '522.77'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 5)
$cell.Value = "'" + '  -0.12%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 4)
$cell.Value = "'" + '139.24'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 5)
$cell.Value = "'" + '  -3.19%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 5)
$cell.Value = "'" + '  +0.24%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 4)
$cell.Value = "'" + '0.563'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 5)
$cell.Value = "'" + '  -1.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 4)
$cell.Value = "'" + '2.599.97'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 5)
$cell.Value = "'" + '  -0.69%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 4)
$cell.Value = "'" + '6.51'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 5)
$cell.Value = "'" + '  -2.00%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 5)
$cell.Value = "'" + '  -0.45%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 5)
$cell.Value = "'" + '  -2.12%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 5)
$cell.Value = "'" + '  +2.57%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'" + '3.047.25'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 5)
$cell.Value = "'" + '  -0.13%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 4)
$cell.Value = "'" + '58.996.38'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 5)
$cell.Value = "'" + '  +1.44%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'" + '20.48'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 5)
$cell.Value = "'" + '  -0.25%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 4)
$cell.Value = "'" + '2.613.31'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 5)
$cell.Value = "'" + '  +1.26%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 5)
$cell.Value = "'" + '  -1.10%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 4)
$cell.Value = "'" + '341.21'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 5)
$cell.Value = "'" + '  +0.53%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 5)
$cell.Value = "'" + '  -1.38%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 4)
$cell.Value = "'" + '10.06'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 5)
$cell.Value = "'" + '  -2.28%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 5)
$cell.Value = "'" + '  +1.01%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 5)
$cell.Value = "'" + '  +0.10%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'" + '66.42'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 5)
$cell.Value = "'" + '  +1.82%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(25, 5)
$cell.Value = "'" + '  +0.73%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 5)
$cell.Value = "'" + '  +0.23%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 4)
$cell.Value = "'" + '0.998'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 5)
$cell.Value = "'" + '  +0.14%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 4)
$cell.Value = "'" + '7.04'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 5)
$cell.Value = "'" + '  +0.38%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 5)
$cell.Value = "'" + '  +0.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(30, 5)
$cell.Value = "'" + '  -3.53%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'" + '5.89'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 5)
$cell.Value = "'" + '  -5.43%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 4)
$cell.Value = "'" + '1.59'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 5)
$cell.Value = "'" + '  -0.16%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 5)
$cell.Value = "'" + '  -0.55%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 5)
$cell.Value = "'" + '  -1.84%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 5)
$cell.Value = "'" + '  -2.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 4)
$cell.Value = "'" + '36.77'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 5)
$cell.Value = "'" + '  +2.04%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 4)
$cell.Value = "'" + '1.46'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 5)
$cell.Value = "'" + '  +0.40%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 5)
$cell.Value = "'" + '  -4.20%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 4)
$cell.Value = "'" + '0.815'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 5)
$cell.Value = "'" + '  -6.53%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 5)
$cell.Value = "'" + '  -0.96%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 4)
$cell.Value = "'" + '0.999'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 5)
$cell.Value = "'" + '  +0.31%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 4)
$cell.Value = "'" + '271.95'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 5)
$cell.Value = "'" + '  -0.61%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 5)
$cell.Value = "'" + '  +0.06%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 4)
$cell.Value = "'" + '10.78'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 5)
$cell.Value = "'" + '  +1.01%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 5)
$cell.Value = "'" + '  -0.94%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 4)
$cell.Value = "'" + '0.0514'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 5)
$cell.Value = "'" + '  -1.78%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(48, 5)
$cell.Value = "'" + '  -2.66%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 4)
$cell.Value = "'" + '1.967.85'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 5)
$cell.Value = "'" + '  -0.34%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'" + '0.0222'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 5)
$cell.Value = "'" + '  -0.28%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 2)
$cell.Value = "'" + 'RenderToken'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 3)
$cell.Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'" + '4.46'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 5)
$cell.Value = "'" + '  -4.06%  '
$cell.Style = 'Normal'
